# Update column G ("K" = strikeouts per game proxy) values for rows 2-10
# with recomputed values (commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 2
    4  = 3
    5  = 4
    6  = 6
    7  = 2
    8  = 3
    9  = 3
    10 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
